$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.475.17'
$ws.Range('E2').Value = '  +2.95%  '
$ws.Range('D3').Value = '2.509.60'
$ws.Range('E3').Value = '  +2.52%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'324.63"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.80%  '
$ws.Range('D6').Value = "'109.83"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.69%  '
$ws.Range('D7').Value = "'0.527"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.72%  '
$ws.Range('D8').Value = "'1.00"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = "'0.543"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('D10').Value = "'39.25"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.74%  '
$ws.Range('D11').Value = "'0.0820"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.70%  '
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').Value = "'18.65"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.71%  '
$ws.Range('E14').Value = '  +2.42%  '
$ws.Range('D15').Value = '2.901.64'
$ws.Range('E15').Value = '  +2.23%  '
$ws.Range('D16').Value = '2.503.61'
$ws.Range('E16').Value = '  +1.72%  '
$ws.Range('E17').Value = '  +2.84%  '
$ws.Range('D18').Value = '47.427.46'
$ws.Range('E18').Value = '  +3.26%  '
$ws.Range('D19').Value = "'12.99"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.59%  '
$ws.Range('E20').Value = '  +4.72%  '
$ws.Range('D21').Value = '0.0₃0949'
$ws.Range('E21').Value = '  +1.41%  '
$ws.Range('E22').Value = '  +11.48%  '
$ws.Range('D23').Value = "'71.01"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('D24').Value = "'250.25"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('E25').Value = '  +3.96%  '
$ws.Range('D26').Value = "'26.26"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.86%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = "'2.29"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.63%  '
$ws.Range('E29').Value = '  +3.84%  '
$ws.Range('D30').Value = "'36.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.49%  '
$ws.Range('E31').Value = '  +4.67%  '
$ws.Range('D32').Value = "'50.33"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.19%  '
$ws.Range('D33').Value = "'19.99"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.95%  '
$ws.Range('D34').Value = "'5.47"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.85%  '
$ws.Range('E35').Value = '  +4.48%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  +6.10%  '
$ws.Range('E38').Value = '  +4.47%  '
$ws.Range('E39').Value = '  +2.85%  '
$ws.Range('E40').Value = '  +1.80%  '
$ws.Range('D41').Value = "'123.18"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.05%  '
$ws.Range('D42').Value = "'2.26"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.84%  '
$ws.Range('D43').Value = "'21.52"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.37%  '
$ws.Range('E44').Value = '  +2.59%  '
$ws.Range('D45').Value = '1.998.70'
$ws.Range('E45').Value = '  +1.99%  '
$ws.Range('E46').Value = '  +4.51%  '
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('E48').Value = '  -2.23%  '
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('E50').Value = '  +8.41%  '
$ws.Range('D51').Value = "'78.79"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.16%  '
